$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1565637.2
$ws.Range("I17").Value = 663
$ws.Range("J17").Value = 1590478.1
$ws.Range("K17").Value = 1989
$ws.Range("L17").Value = 4771434.300000001
$ws.Range("M17").Value = -1821
$ws.Range("N17").Value = -4771770.300000001
$ws.Range("H88").Value = 1277.1578
$ws.Range("I88").Value = 917.625
$ws.Range("J88").Value = 1538.6364
$ws.Range("K88").Value = 917.625
$ws.Range("L88").Value = 1538.6364
$ws.Range("M88").Value = -511.625
$ws.Range("N88").Value = -2350.6364
$ws.Range("H91").Value = 1277.1578
$ws.Range("I91").Value = 917.625
$ws.Range("J91").Value = 1538.6364
$ws.Range("K91").Value = 917.625
$ws.Range("L91").Value = 1538.6364
$ws.Range("M91").Value = 486.375
$ws.Range("N91").Value = -4346.6364
$ws.Range("H127").Value = 1805.8
$ws.Range("I127").Value = 846.3333
$ws.Range("K127").Value = 2538.9999
$ws.Range("M127").Value = 2421.0001
$ws.Range("H129").Value = 271063.94
$ws.Range("J129").Value = 294951.94
$ws.Range("L129").Value = 884855.8200000001
$ws.Range("N129").Value = -894855.8200000001
$ws.Range("H131").Value = 4650
$ws.Range("I131").Value = 4500
$ws.Range("J131").Value = 4800
$ws.Range("K131").Value = 13500
$ws.Range("L131").Value = 14400
$ws.Range("M131").Value = -8460
$ws.Range("N131").Value = -24480
$ws.Range("H132").Value = 3129.3125
$ws.Range("I132").Value = 3274.077
$ws.Range("K132").Value = 9822.231
$ws.Range("M132").Value = -7292.231
$ws.Range("H138").Value = 1426.0513
$ws.Range("I138").Value = 521.56525
$ws.Range("K138").Value = 1564.69575
$ws.Range("M138").Value = 3575.30425

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 10000
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10368
$ws.Range("H32").Value = 22675.06
$ws.Range("I32").Value = 23324.234
$ws.Range("J32").Value = 12504.667
$ws.Range("K32").Value = 23324.234
$ws.Range("L32").Value = 12504.667
$ws.Range("M32").Value = -23037.234
$ws.Range("N32").Value = -13078.667
$ws.Range("H61").Value = 2754.28
$ws.Range("I61").Value = 2234
$ws.Range("K61").Value = 2234
$ws.Range("M61").Value = -2022
$ws.Range("H74").Value = 52632356
$ws.Range("I74").Value = 100000380
$ws.Range("J74").Value = 1209.2222
$ws.Range("K74").Value = 100000380
$ws.Range("L74").Value = 1209.2222
$ws.Range("M74").Value = -99999506
$ws.Range("N74").Value = -2957.2222
$ws.Range("H77").Value = 52632356
$ws.Range("I77").Value = 100000380
$ws.Range("J77").Value = 1209.2222
$ws.Range("K77").Value = 500001900
$ws.Range("L77").Value = 6046.111
$ws.Range("M77").Value = -499997532
$ws.Range("N77").Value = -14782.111
$ws.Range("H122").Value = 1290.762
$ws.Range("I122").Value = 1331.8948
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 3995.6844
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -1545.6844
$ws.Range("N122").Value = -7600
$ws.Range("H136").Value = 2754.28
$ws.Range("I136").Value = 2234
$ws.Range("K136").Value = 6702
$ws.Range("M136").Value = -4152

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1402.9333
$ws.Range("I80").Value = 1322.7693
$ws.Range("J80").Value = 1464.2354
$ws.Range("K80").Value = 1322.7693
$ws.Range("L80").Value = 1464.2354
$ws.Range("M80").Value = -324.7692999999999
$ws.Range("N80").Value = -3460.2354
$ws.Range("H83").Value = 1402.9333
$ws.Range("I83").Value = 1322.7693
$ws.Range("J83").Value = 1464.2354
$ws.Range("K83").Value = 6613.8465
$ws.Range("L83").Value = 7321.177
$ws.Range("M83").Value = -1621.8465
$ws.Range("N83").Value = -17305.177
$ws.Range("H134").Value = 30607.217
$ws.Range("I134").Value = 38497.93
$ws.Range("K134").Value = 115493.79
$ws.Range("M134").Value = -112958.79

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1260
$ws.Range("I16").Value = 1390
$ws.Range("K16").Value = 1390
$ws.Range("M16").Value = -1103
$ws.Range("H22").Value = 207.14285
$ws.Range("I22").Value = 187.5
$ws.Range("J22").Value = 233.33333
$ws.Range("K22").Value = 187.5
$ws.Range("L22").Value = 233.33333
$ws.Range("M22").Value = 162.5
$ws.Range("N22").Value = -933.3333299999999
$ws.Range("H31").Value = 9332.046
$ws.Range("I31").Value = 11401.936
$ws.Range("J31").Value = 4396.154
$ws.Range("K31").Value = 11401.936
$ws.Range("L31").Value = 4396.154
$ws.Range("M31").Value = -11106.936
$ws.Range("N31").Value = -4986.154
$ws.Range("H34").Value = 9332.046
$ws.Range("I34").Value = 11401.936
$ws.Range("J34").Value = 4396.154
$ws.Range("K34").Value = 11401.936
$ws.Range("L34").Value = 4396.154
$ws.Range("M34").Value = -11199.936
$ws.Range("N34").Value = -4800.154
$ws.Range("H58").Value = 15353.028
$ws.Range("I58").Value = 1094.6
$ws.Range("K58").Value = 1094.6
$ws.Range("M58").Value = -891.5999999999999
$ws.Range("H113").Value = 1260
$ws.Range("I113").Value = 1390
$ws.Range("K113").Value = 1390
$ws.Range("M113").Value = 780
$ws.Range("H132").Value = 23633.72
$ws.Range("I132").Value = 29033.475
$ws.Range("J132").Value = 6534.5
$ws.Range("K132").Value = 87100.42499999999
$ws.Range("L132").Value = 19603.5
$ws.Range("M132").Value = -84570.42499999999
$ws.Range("N132").Value = -24663.5
$ws.Range("H134").Value = 1244.8667
$ws.Range("I134").Value = 926.34784
$ws.Range("K134").Value = 2779.04352
$ws.Range("M134").Value = -244.0435200000002
$ws.Range("H136").Value = 15353.028
$ws.Range("I136").Value = 1094.6
$ws.Range("K136").Value = 3283.8
$ws.Range("M136").Value = -733.7999999999997

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 676.9
$ws.Range("I40").Value = 96.125
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 384.5
$ws.Range("L40").Value = 12000
$ws.Range("M40").Value = -315.5
$ws.Range("N40").Value = -12138
$ws.Range("H131").Value = 760.33
$ws.Range("J131").Value = 771.9367999999999
$ws.Range("L131").Value = 2315.8104
$ws.Range("N131").Value = -12395.8104

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 588.58826
$ws.Range("I97").Value = 502.15384
$ws.Range("J97").Value = 869.5
$ws.Range("K97").Value = 502.15384
$ws.Range("L97").Value = 869.5
$ws.Range("M97").Value = -6.153840000000002
$ws.Range("N97").Value = -1861.5
$ws.Range("H113").Value = 3230.7693
$ws.Range("I113").Value = 2555.5557
$ws.Range("J113").Value = 4750
$ws.Range("K113").Value = 2555.5557
$ws.Range("L113").Value = 4750
$ws.Range("M113").Value = -385.5556999999999
$ws.Range("N113").Value = -9090

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1937.5625
$ws.Range("I22").Value = 1700.1111
$ws.Range("J22").Value = 2242.8572
$ws.Range("K22").Value = 1700.1111
$ws.Range("L22").Value = 2242.8572
$ws.Range("M22").Value = -1405.1111
$ws.Range("N22").Value = -2832.8572
$ws.Range("H27").Value = 1937.5625
$ws.Range("I27").Value = 1700.1111
$ws.Range("J27").Value = 2242.8572
$ws.Range("K27").Value = 1700.1111
$ws.Range("L27").Value = 2242.8572
$ws.Range("M27").Value = -1593.1111
$ws.Range("N27").Value = -2456.8572
$ws.Range("H117").Value = 39990
$ws.Range("J117").Value = 39990
$ws.Range("L117").Value = 39990
$ws.Range("N117").Value = -49168
$ws.Range("H132").Value = 2523.3333
$ws.Range("I132").Value = 1666.1333
$ws.Range("J132").Value = 4666.3335
$ws.Range("K132").Value = 4998.3999
$ws.Range("L132").Value = 13999.0005
$ws.Range("M132").Value = -2468.3999
$ws.Range("N132").Value = -19059.0005
$ws.Range("H136").Value = 25571.047
$ws.Range("I136").Value = 36856.57
$ws.Range("K136").Value = 110569.71
$ws.Range("M136").Value = -108019.71

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 40342
$ws.Range("J118").Value = 40342
$ws.Range("L118").Value = 40342
$ws.Range("N118").Value = -43656
$ws.Range("H132").Value = 1917.8182
$ws.Range("I132").Value = 1312.3077
$ws.Range("J132").Value = 2792.4443
$ws.Range("K132").Value = 3936.9231
$ws.Range("L132").Value = 8377.332900000001
$ws.Range("M132").Value = -1406.9231
$ws.Range("N132").Value = -13437.3329

Write-Host "Applied 215 cell updates across 8 sheets."